# Insert a new weekly price record as row 131 (Albahaca / Feria Lagunitas de
# Puerto Montt). Inserting the row pushes the existing rows 131-166 down to
# 132-167, preserving all of their data and formatting, and the sheet's used
# range grows from A1:R166 to A1:R167.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(131).Insert()

$ws.Cells.Item(131, 1).Value = 4
$ws.Cells.Item(131, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(131, 3).Value = "Los Lagos"
$ws.Cells.Item(131, 4).Value = 44964
$ws.Cells.Item(131, 5).Value = 10
$ws.Cells.Item(131, 6).Value = 100112052
$ws.Cells.Item(131, 7).Value = "Albahaca"
$ws.Cells.Item(131, 8).Value = "Sin especificar"
$ws.Cells.Item(131, 9).Value = "Primera"
$ws.Cells.Item(131, 10).Value = 120
$ws.Cells.Item(131, 11).Value = 6000
$ws.Cells.Item(131, 12).Value = 6000
$ws.Cells.Item(131, 13).Value = 6000
$ws.Cells.Item(131, 14).Value = "`$/docena de matas"
$ws.Cells.Item(131, 15).Value = "Región Metropolitana"
$ws.Cells.Item(131, 16).Value = 1000
$ws.Cells.Item(131, 17).Value = 6
$ws.Cells.Item(131, 18).Value = "Hortaliza"

# Keep the date column using the same date number format as the rest of the
# column (style carries over from the Insert, but set it explicitly so it's
# not left as a plain number if the insert didn't propagate formatting).
$ws.Cells.Item(131, 4).NumberFormat = $ws.Cells.Item(132, 4).NumberFormat
